$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.139.98"
$ws.Range("E2").Value = "  -0.22%  "
$ws.Range("D3").Value = "2.648.13"
$ws.Range("E3").Value = "  +2.75%  "
$ws.Range("E4").Value = "  +0.33%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "609.71"
$ws.Range("E5").Value = "  +3.70%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.21"
$ws.Range("E6").Value = "  -1.15%  "
$ws.Range("E7").Value = "  +0.21%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.585"
$ws.Range("E8").Value = "  -1.04%  "
$ws.Range("D9").Value = "2.644.16"
$ws.Range("E9").Value = "  +2.66%  "
$ws.Range("E10").Value = "  +0.00%  "
$ws.Range("E11").Value = "  +0.04%  "
$ws.Range("E12").Value = "  +0.38%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.360"
$ws.Range("E13").Value = "  +2.51%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.18"
$ws.Range("E14").Value = "  -0.23%  "
$ws.Range("D15").Value = "3.130.40"
$ws.Range("E15").Value = "  +3.05%  "
$ws.Range("D16").Value = "63.079.98"
$ws.Range("E16").Value = "  -0.19%  "
$ws.Range("E17").Value = "  -1.37%  "
$ws.Range("D18").Value = "2.647.74"
$ws.Range("E18").Value = "  +2.75%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.40"
$ws.Range("E19").Value = "  +2.86%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "341.06"
$ws.Range("E20").Value = "  -0.18%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.40"
$ws.Range("E21").Value = "  +1.35%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.84"
$ws.Range("E22").Value = "  +2.96%  "
$ws.Range("E23").Value = "  -0.16%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "66.88"
$ws.Range("E24").Value = "  -1.63%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.64"
$ws.Range("E25").Value = "  +1.42%  "
$ws.Range("E26").Value = "  -1.95%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.61"
$ws.Range("E27").Value = "  +4.44%  "
$ws.Range("E28").Value = "  -0.72%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "545.20"
$ws.Range("E29").Value = "  +15.27%  "
$ws.Range("E30").Value = "  +0.08%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.78"
$ws.Range("E31").Value = "  -2.15%  "
$ws.Range("E32").Value = "  +5.43%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.79"
$ws.Range("E33").Value = "  +5.94%  "
$ws.Range("D34").Value = "0.0₃0804"
$ws.Range("E34").Value = "  +0.05%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "173.13"
$ws.Range("E35").Value = "  -1.80%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.13"
$ws.Range("E36").Value = "  +12.79%  "
$ws.Range("E37").Value = "  +1.34%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.996"
$ws.Range("E38").Value = "  -0.50%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.08"
$ws.Range("E39").Value = "  +1.01%  "
$ws.Range("E40").Value = "  +8.39%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "175.03"
$ws.Range("E41").Value = "  +10.99%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  +0.14%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.73"
$ws.Range("E43").Value = "  +0.82%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "22.09"
$ws.Range("E44").Value = "  +3.58%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0573"
$ws.Range("E45").Value = "  +6.02%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.630"
$ws.Range("E46").Value = "  -0.72%  "
$ws.Range("E47").Value = "  +1.11%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0960"
$ws.Range("E48").Value = "  -0.37%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "18.60"
$ws.Range("E49").Value = "  +2.45%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.73"
$ws.Range("E50").Value = "  +3.40%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "11.32"
$ws.Range("E51").Value = "  -0.67%  "
